# Auto-generated edit script: updates market-price-derived profit columns (H-N)
# across the 8 class sheets, per the scheduled runner's refreshed pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 12602593
$ws.Range("I116").Value = 14645043
$ws.Range("K116").Value = 14645043
$ws.Range("M116").Value = -14641601
$ws.Range("H129").Value = 1523.5
$ws.Range("I129").Value = 1064.6666
$ws.Range("J129").Value = 2900
$ws.Range("K129").Value = 3193.9998
$ws.Range("L129").Value = 8700
$ws.Range("M129").Value = 1806.0002
$ws.Range("N129").Value = -18700
$ws.Range("H132").Value = 148995.7
$ws.Range("I132").Value = 217888.62
$ws.Range("K132").Value = 653665.86
$ws.Range("M132").Value = -651135.86
$ws.Range("H138").Value = 2881.899
$ws.Range("I138").Value = 1441.5306
$ws.Range("J138").Value = 4646.35
$ws.Range("K138").Value = 4324.5918
$ws.Range("L138").Value = 13939.05
$ws.Range("M138").Value = 815.4081999999999
$ws.Range("N138").Value = -24219.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11294.965
$ws.Range("I32").Value = 11498.62
$ws.Range("J32").Value = 10331
$ws.Range("K32").Value = 11498.62
$ws.Range("L32").Value = 10331
$ws.Range("M32").Value = -11211.62
$ws.Range("N32").Value = -10905
$ws.Range("H63").Value = 6248.75
$ws.Range("I63").Value = 5000
$ws.Range("K63").Value = 5000
$ws.Range("M63").Value = -4314
$ws.Range("H66").Value = 6248.75
$ws.Range("I66").Value = 5000
$ws.Range("K66").Value = 25000
$ws.Range("M66").Value = -21568
$ws.Range("H107").Value = 105000
$ws.Range("J107").Value = 105000
$ws.Range("L107").Value = 105000
$ws.Range("N107").Value = -112680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H76").Value = 66607
$ws.Range("J76").Value = 66607
$ws.Range("L76").Value = 66607
$ws.Range("N76").Value = -67237
$ws.Range("H79").Value = 66607
$ws.Range("J79").Value = 66607
$ws.Range("L79").Value = 66607
$ws.Range("N79").Value = -68791
$ws.Range("H105").Value = 44119776
$ws.Range("I105").Value = 44119776
$ws.Range("K105").Value = 44119776
$ws.Range("M105").Value = -44118029
$ws.Range("H134").Value = 2229.8367
$ws.Range("I134").Value = 2048.3157
$ws.Range("K134").Value = 6144.9471
$ws.Range("M134").Value = -3609.9471

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1173.1818
$ws.Range("I22").Value = 687
$ws.Range("K22").Value = 687
$ws.Range("M22").Value = -337
$ws.Range("H31").Value = 20835950
$ws.Range("I31").Value = 26317906
$ws.Range("K31").Value = 26317906
$ws.Range("M31").Value = -26317611
$ws.Range("H34").Value = 20835950
$ws.Range("I34").Value = 26317906
$ws.Range("K34").Value = 26317906
$ws.Range("M34").Value = -26317704
$ws.Range("H74").Value = 79999.39999999999
$ws.Range("J74").Value = 79999.39999999999
$ws.Range("L74").Value = 79999.39999999999
$ws.Range("N74").Value = -81747.39999999999
$ws.Range("H77").Value = 79999.39999999999
$ws.Range("J77").Value = 79999.39999999999
$ws.Range("L77").Value = 239998.2
$ws.Range("N77").Value = -248734.2
$ws.Range("H97").Value = 31495.834
$ws.Range("J97").Value = 31495.834
$ws.Range("L97").Value = 31495.834
$ws.Range("N97").Value = -33477.834
$ws.Range("H107").Value = 957599.6
$ws.Range("I107").Value = 1212699
$ws.Range("K107").Value = 1212699
$ws.Range("M107").Value = -1210779
$ws.Range("H111").Value = 87993
$ws.Range("J111").Value = 87993
$ws.Range("L111").Value = 87993
$ws.Range("N111").Value = -96173
$ws.Range("H132").Value = 30314252
$ws.Range("I132").Value = 41674900
$ws.Range("K132").Value = 125024700
$ws.Range("M132").Value = -125022170
$ws.Range("H141").Value = 168414.33
$ws.Range("J141").Value = 216896.75
$ws.Range("L141").Value = 216896.75
$ws.Range("N141").Value = -227256.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 55563316
$ws.Range("I56").Value = 55563316
$ws.Range("K56").Value = 55563316
$ws.Range("M56").Value = -55562786
$ws.Range("H121").Value = 1551.2
$ws.Range("J121").Value = 1551.2
$ws.Range("L121").Value = 4653.6
$ws.Range("N121").Value = -7273.6
$ws.Range("H122").Value = 1350
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 10800
$ws.Range("M122").Value = -8350
$ws.Range("H131").Value = 22545220
$ws.Range("I131").Value = 33334254
$ws.Range("J131").Value = 21046742
$ws.Range("K131").Value = 100002762
$ws.Range("L131").Value = 63140226
$ws.Range("M131").Value = -99997722
$ws.Range("N131").Value = -63150306

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3442.3845
$ws.Range("J113").Value = 3510.25
$ws.Range("L113").Value = 3510.25
$ws.Range("N113").Value = -7850.25
$ws.Range("H122").Value = 320006.6
$ws.Range("I122").Value = 734897.2
$ws.Range("K122").Value = 2204691.6
$ws.Range("M122").Value = -2202241.6
$ws.Range("H123").Value = 56406.4
$ws.Range("J123").Value = 56406.4
$ws.Range("L123").Value = 56406.4
$ws.Range("N123").Value = -61306.4
$ws.Range("H126").Value = 4025.4583
$ws.Range("I126").Value = 2278.1765
$ws.Range("J126").Value = 8268.857
$ws.Range("K126").Value = 6834.529500000001
$ws.Range("L126").Value = 24806.571
$ws.Range("M126").Value = -4364.529500000001
$ws.Range("N126").Value = -29746.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2602.4575
$ws.Range("I132").Value = 2408.6882
$ws.Range("K132").Value = 7226.0646
$ws.Range("M132").Value = -4696.0646
$ws.Range("H136").Value = 5193.0586
$ws.Range("I136").Value = 2745
$ws.Range("K136").Value = 8235
$ws.Range("M136").Value = -5685
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 11499.615
$ws.Range("I18").Value = 9333.333000000001
$ws.Range("K18").Value = 9333.333000000001
$ws.Range("M18").Value = -9160.333000000001
$ws.Range("H70").Value = 44995
$ws.Range("I70").Value = 44990
$ws.Range("K70").Value = 44990
$ws.Range("M70").Value = -44675
$ws.Range("H73").Value = 44995
$ws.Range("I73").Value = 44990
$ws.Range("K73").Value = 44990
$ws.Range("M73").Value = -43898
$ws.Range("H81").Value = 2103876.8
$ws.Range("I81").Value = 5213208
$ws.Range("J81").Value = 30989.166
$ws.Range("K81").Value = 10426416
$ws.Range("L81").Value = 61978.332
$ws.Range("M81").Value = -10425355
$ws.Range("N81").Value = -64100.332
$ws.Range("H84").Value = 2103876.8
$ws.Range("I84").Value = 5213208
$ws.Range("J84").Value = 30989.166
$ws.Range("K84").Value = 52132080
$ws.Range("L84").Value = 309891.66
$ws.Range("M84").Value = -52126776
$ws.Range("N84").Value = -320499.66
$ws.Range("H122").Value = 4115.6313
$ws.Range("I122").Value = 3261.5557
$ws.Range("J122").Value = 6212
$ws.Range("K122").Value = 9784.667099999999
$ws.Range("L122").Value = 18636
$ws.Range("M122").Value = -7334.667099999999
$ws.Range("N122").Value = -23536
$ws.Range("H132").Value = 2551.761
$ws.Range("I132").Value = 1119.5676
$ws.Range("J132").Value = 8439.666999999999
$ws.Range("K132").Value = 3358.7028
$ws.Range("L132").Value = 25319.001
$ws.Range("M132").Value = -828.7028
$ws.Range("N132").Value = -30379.001
$ws.Range("H136").Value = 6159.1123
$ws.Range("I136").Value = 3633.9153
$ws.Range("J136").Value = 9979.281999999999
$ws.Range("K136").Value = 10901.7459
$ws.Range("L136").Value = 29937.846
$ws.Range("M136").Value = -8351.7459
$ws.Range("N136").Value = -35037.846
